$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the three rows (identified by their "Materia" text) that were removed
# from the dataset. Row 9 = COMPUTO FLEXIBLE (SOFTCOMPUTING); the former rows 23
# and 24 = PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION and
# PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS.
$ws.Range("A9:D9").EntireRow.Delete() | Out-Null
$ws.Range("A22:D22").EntireRow.Delete() | Out-Null
$ws.Range("A22:D22").EntireRow.Delete() | Out-Null
